# Updated cryptos list values (price, volume %, and re-ranked rows 34-36)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'" + '22.002.11'
$ws.Range("E2").Value = '  -1.87%  '
$ws.Range("D3").Formula = "'" + '1.553.94'
$ws.Range("E3").Value = '  -0.97%  '
$ws.Range("D4").Formula = "'" + '1.000'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").Formula = "'" + '286.90'
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").Formula = "'" + '0.3801'
$ws.Range("E7").Value = '  +2.96%  '
$ws.Range("D8").Formula = "'" + '0.3245'
$ws.Range("E8").Value = '  -2.14%  '
$ws.Range("D9").Formula = "'" + '41.57'
$ws.Range("E9").Value = '  -12.09%  '
$ws.Range("D10").Formula = "'" + '1.120'
$ws.Range("E10").Value = '  -3.23%  '
$ws.Range("D11").Formula = "'" + '0.07316'
$ws.Range("E11").Value = '  -2.17%  '
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").Formula = "'" + '19.37'
$ws.Range("E13").Value = '  -6.43%  '
$ws.Range("D14").Formula = "'" + '5.719'
$ws.Range("E14").Value = '  -3.49%  '
$ws.Range("D15").Formula = "'" + '6.807'
$ws.Range("E15").Value = '  -0.99%  '
$ws.Range("D16").Formula = "'" + '1.558.68'
$ws.Range("E16").Value = '  -0.04%  '
$ws.Range("D17").Formula = "'" + '0.00001095'
$ws.Range("E17").Value = '  -1.44%  '
$ws.Range("D18").Formula = "'" + '0.06627'
$ws.Range("E18").Value = '  -1.46%  '
$ws.Range("D19").Formula = "'" + '85.09'
$ws.Range("E19").Value = '  -2.89%  '
$ws.Range("D20").Formula = "'" + '6.424'
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("D21").Formula = "'" + '0.9996'
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("E22").Value = '  -2.88%  '
$ws.Range("D23").Formula = "'" + '11.48'
$ws.Range("E23").Value = '  -3.79%  '
$ws.Range("D24").Formula = "'" + '22.022.94'
$ws.Range("E24").Value = '  -1.76%  '
$ws.Range("D25").Formula = "'" + '2.296'
$ws.Range("E25").Value = '  -3.20%  '
$ws.Range("D26").Formula = "'" + '2.521'
$ws.Range("E26").Value = '  -3.57%  '
$ws.Range("D27").Formula = "'" + '148.70'
$ws.Range("E27").Value = '  -1.40%  '
$ws.Range("D28").Formula = "'" + '18.86'
$ws.Range("E28").Value = '  -3.37%  '
$ws.Range("D29").Formula = "'" + '4.843'
$ws.Range("E29").Value = '  -1.79%  '
$ws.Range("D30").Formula = "'" + '1.732.23'
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("D31").Formula = "'" + '120.52'
$ws.Range("E31").Value = '  -3.14%  '
$ws.Range("D32").Formula = "'" + '1.099'
$ws.Range("E32").Value = '  +1.53%  '
$ws.Range("D33").Formula = "'" + '5.899'
$ws.Range("E33").Value = '  -2.80%  '
$ws.Range("B34").Formula = "'" + 'FraxShare'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D34").Formula = "'" + '9.281'
$ws.Range("E34").Value = '  -5.65%  '
$ws.Range("B35").Formula = "'" + 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").Formula = "'" + '0.08151'
$ws.Range("E35").Value = '  -2.02%  '
$ws.Range("B36").Formula = "'" + 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Formula = "'" + '1.651'
$ws.Range("E36").Value = '  -16.78%  '
$ws.Range("D37").Formula = "'" + '0.06201'
$ws.Range("E37").Value = '  -2.78%  '
$ws.Range("D38").Formula = "'" + '5.245'
$ws.Range("E38").Value = '  -1.21%  '
$ws.Range("D39").Formula = "'" + '0.02290'
$ws.Range("E39").Value = '  -5.84%  '
$ws.Range("E40").Value = '  -4.71%  '
$ws.Range("D41").Formula = "'" + '1.221'
$ws.Range("E41").Value = '  -5.85%  '
$ws.Range("E42").Value = '  -4.27%  '
$ws.Range("D43").Formula = "'" + '1.0000'
$ws.Range("D44").Formula = "'" + '0.5938'
$ws.Range("E44").Value = '  -3.78%  '
$ws.Range("D45").Formula = "'" + '13.52'
$ws.Range("E45").Value = '  -3.40%  '
$ws.Range("D46").Formula = "'" + '3.724'
$ws.Range("E46").Value = '  -1.07%  '
$ws.Range("D47").Formula = "'" + '0.5735'
$ws.Range("E47").Value = '  -4.30%  '
$ws.Range("D48").Formula = "'" + '1.933'
$ws.Range("E48").Value = '  -5.02%  '
$ws.Range("D49").Formula = "'" + '119.48'
$ws.Range("E49").Value = '  -4.33%  '
$ws.Range("D50").Formula = "'" + '1.156'
$ws.Range("E50").Value = '  -3.39%  '
$ws.Range("D51").Formula = "'" + '0.06872'
$ws.Range("E51").Value = '  -4.38%  '
